$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "Hoja1" -> "recibos"
$ws.Name = "recibos"

# Update the week label (shared string used by B9; H9/B27/H27/B43 are
# formulas that reference it and recompute automatically).
$ws.Range("B9").Value = "SEMANA  30  DEL    25      Al   31   DE   JULIO          2022"

# Bonus/extra amount for the second block went from 560 to 1400
# (K24 = SUM(K21:K23) recalculates automatically).
$ws.Range("K21").Value = 1400

# Update the window's active cell / selection for the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I39").Select()
